$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Binance -> binance)
$ws.Name = "binance"

# Replace the old single data row with a header row of 6 columns.
$headers = @("Time", "Exchange", "Arbitrage Direction", "Cryptocurrency Pairs", "Initial Investment", "Profit/Loss")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# F1 is a brand-new cell with no prior style - copy the existing bold/border/
# centered style (already applied to B1:E1) onto it so all headers match.
$copyResult = $ws.Range("B1").Copy()
$pasteResult = $ws.Range("F1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "Profit/Loss"

# Time column gets an explicit time number format.
$ws.Cells.Item(1, 1).NumberFormat = "h:mm:ss"

# Column widths for B..F.
$ws.Columns.Item(2).ColumnWidth = 13.58333333333333
$ws.Columns.Item(3).ColumnWidth = 25.416666666666664
$ws.Columns.Item(4).ColumnWidth = 33.74999999999999
$ws.Columns.Item(5).ColumnWidth = 19.08333333333333
$ws.Columns.Item(6).ColumnWidth = 13.249999999999998

# Match the saved selection/active cell.
$selResult = $ws.Range("C6").Select()
